$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same formatting as the other
# header cells (e.g. G1) by copying/pasting formats, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add new value in H2 for the existing data row
$ws.Range("H2").Value = 1
